$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("other")

# Insert a new blank row at 65 -- this shifts the old row 65 ("**" terminator)
# down to row 66 without touching its formatting (it keeps no explicit style).
$ws.Rows.Item(65).Insert()

# Pre-create cells A65:D65 as "existing" cells with no baked style by copying
# format/type from the style-less cells directly above (row 64), then overwrite
# their value. Because the cells already exist, Excel doesn't fall back to the
# column's default style when the value is set.
$ws.Cells.Item(64, 1).Copy($ws.Cells.Item(65, 1))
$ws.Cells.Item(65, 1).Value = "28/03/2018"

$ws.Cells.Item(64, 2).Copy($ws.Cells.Item(65, 2))
# B65 keeps the same equipment number as the row above, so no value change needed.

$ws.Cells.Item(64, 3).Copy($ws.Cells.Item(65, 3))
# C65 keeps the same value as the row above, so no value change needed.

$ws.Cells.Item(64, 4).Copy($ws.Cells.Item(65, 4))
$ws.Cells.Item(65, 4).Value = "Не відповідне скручення проводів"

$ws.Cells.Item(64, 5).Copy($ws.Cells.Item(65, 5))

# E65 must hold the text "50" (not the number 50). Assigning a purely numeric
# looking string through .Value always gets auto-converted to a number, even
# on a pre-existing text cell, so build it via a throw-away formula cell that
# evaluates to the text "50", copy *that*, and paste-special only the value
# onto the (already existing) E65 cell. Then remove every trace of the helper.
$helperRow = 500
$ws.Cells.Item($helperRow, 1).Formula = '="50"'
$ws.Cells.Item($helperRow, 1).Copy()
$ws.Cells.Item(65, 5).PasteSpecial(-4163) # xlPasteValues
$ws.Rows.Item($helperRow).Delete()
